# Generate Report for Handoff
# Inserts a new row for file "1e4106c6-b7f9-4752-aaf1-1c3084888f01.md"
# (status "Ready for handoff") immediately before the existing
# "556b5bf8-5d8c-45d3-9293-87717bb47a62.md" row, on all three sheets:
#   Overview (sheet1), zh-cn (sheet2), de-de (sheet3)

$wb = $excel.ActiveWorkbook

$newFileBase = "1e4106c6-b7f9-4752-aaf1-1c3084888f01"
$newMdName   = "$newFileBase.md"
$statusReady = "Ready for handoff"

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(6).Insert()

$wsOverview.Range("A6").Value = $newMdName
$wsOverview.Range("B6").Value = $statusReady
$wsOverview.Range("C6").Value = $statusReady

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$newMdName",
    "",
    "",
    $newMdName
)

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows.Item(6).Insert()

$zhXlfName = "$newFileBase.02ca2d424d31281705f41948a7da8e3f245bc385.zh-cn.xlf"

$wsZhCn.Range("A6").Value = $newMdName
$wsZhCn.Range("B6").Value = $statusReady
$wsZhCn.Range("C6").Value = $zhXlfName
$wsZhCn.Range("D6").Value = "2016-03-04 01:23:55"
$wsZhCn.Range("G6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H6").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$newMdName",
    "",
    "",
    $newMdName
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C6"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName",
    "",
    "",
    $zhXlfName
)

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows.Item(6).Insert()

$deXlfName = "$newFileBase.02ca2d424d31281705f41948a7da8e3f245bc385.de-de.xlf"

$wsDeDe.Range("A6").Value = $newMdName
$wsDeDe.Range("B6").Value = $statusReady
$wsDeDe.Range("C6").Value = $deXlfName
$wsDeDe.Range("D6").Value = "2016-03-04 01:24:10"
$wsDeDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H6").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A6"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$newMdName",
    "",
    "",
    $newMdName
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C6"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName",
    "",
    "",
    $deXlfName
)
